# Edit script: "Harmonizing Nature's Melody" (music/math essay) ->
# "Unveiling the Symphony of Life" (biology essay)

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Harmonizing Nature's Melody: Music, Math, and the Aesthetics of Sound" `
             "Unveiling the Symphony of Life: An Exploration of Biology"

# --- Byline / author name ---
Replace-Text " Aurelio De Lira" " Sarah Morgan"

# --- Email address (local part first, then domain); "Aurelio" standalone remains
#     after the byline replacement above, so this only touches the email run. ---
Replace-Text "Aurelio" "sarah"
Replace-Text "DeLira@musicconservatory" "morgan@school"

# --- Intro paragraph, sentence by sentence ---
Replace-Text "Within the vast tapestry of human existence, music occupies a realm of its own, weaving threads of emotion, intellect, and transcendence" `
             "Biology, the study of life, is a captivating and complex field that encompasses the intricate workings of living organisms"

Replace-Text " Its allure lies not only in its ability to stir the soul but also in its profound connection to the intricate web of mathematics and the aesthetics of sound" `
             " From the smallest microorganisms to the grandest whales, biology unveils the harmonious symphony of life and the profound interconnectedness of all living things"

Replace-Text " In this exploration, we delve into the harmonious interplay of music, math, and aesthetics, uncovering the hidden patterns that lend beauty and meaning to the melodies that grace our ears" `
             " As we delve into this realm of knowledge, we embark on a journey of discovery, unraveling the mysteries of life and gaining a deeper understanding of ourselves and the natural world"

Replace-Text "From the subtle vibrations of a plucked string to the intricate harmonies of a symphony, music is governed by mathematical principles that create a foundation of order and symmetry" `
             "Biology unveils the intricate mechanisms underlying the functioning of organisms, revealing the symphony between cells, tissues, organs, and systems"

Replace-Text " These mathematical underpinnings, such as ratios, proportions, and harmonic progressions, serve as the structural backbone of musical compositions, providing a framework for melodic development and variation" `
             " Through meticulous observation and experimentation, biologists have unlocked the secrets of cellular respiration, DNA replication, and protein synthesis--the fundamental processes that sustain life"

Replace-Text " It is this mathematical order that enables musicians to craft cohesive pieces that resonate with our sense of equilibrium and balance" `
             " These discoveries have revolutionized medicine, leading to innovative treatments and therapies"

Replace-Text "Moreover, music's aesthetic appeal stems from its ability to evoke emotions and convey narratives" `
             "Biology not only probes the inner workings of organisms but also delves into the interactions between different species and the delicate balance of ecosystems"

Replace-Text " Whether it's the melancholic strains of a minor key or the uplifting melodies of a major chord, music possesses an uncanny ability to tap into our deepest emotions, triggering memories, and creating a sense of connection with others" `
             " Studies of symbiotic relationships, food chains, and biogeochemical cycles underscore the interconnectedness of life and the importance of preserving biodiversity"

# --- Remove the rest of the old essay body (from "This emotional resonance..."
#     through "...hearts and minds of people across cultures, time, and space")
#     and replace it with a single closing sentence, leaving the paragraph's
#     trailing period run untouched. ---
$bigPara = $d.Paragraphs(5).Range
$bigText = $bigPara.Text
$cutStart = $bigText.IndexOf(" This emotional resonance")
$cutEnd = $bigText.IndexOf(" It is through this interplay that music transcends its physical form, becoming a universal language capable of speaking to the hearts and minds of people across cultures, time, and space")
$cutEnd = $cutEnd + " It is through this interplay that music transcends its physical form, becoming a universal language capable of speaking to the hearts and minds of people across cultures, time, and space".Length
$rStart = $bigPara.Start + $cutStart
$rEnd = $bigPara.Start + $cutEnd
$cutRange = $d.Range($rStart, $rEnd)
$cutRange.Delete() | Out-Null
$insertionPoint = $d.Range($rStart, $rStart)
$insertionPoint.InsertBefore(" Biology equips us with the knowledge to address pressing environmental challenges and promote sustainable practices") | Out-Null

# --- Summary paragraph ---
Replace-Text "In this exploration of the relationship between music, mathematics, and aesthetics, we have uncovered the intricate web of principles that lend beauty and meaning to the melodies that grace our ears" `
             "Biology is a captivating journey of discovery, unraveling the complexities of life and the interconnectedness of living organisms"

Replace-Text " Music's mathematical foundation provides a framework for creating cohesive compositions, while aesthetic principles shape its emotional impact and structure" `
             " Through the study of cells, organisms, and ecosystems, we gain profound insights into the mechanisms underlying life and the intricate web of relationships that sustain our planet"

Replace-Text " This harmonious interplay enables music to transcend its physical form, becoming a universal language capable of speaking to the deepest recesses of the human experience" `
             " Biology empowers us to address global challenges, understand our place in the natural world, and appreciate the mesmerizing symphony of life"

# --- Add a new, completely empty paragraph at the very end of the document
#     body (after the Summary paragraph, before the section properties). Using
#     raw XML insertion (rather than InsertParagraphAfter) avoids carrying the
#     previous run's formatting into a stray empty run on the new mark, giving
#     a bare <w:p/> like the target. ---
$endRange = $d.Content
$endRange.Collapse(0) | Out-Null
$endRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>") | Out-Null
